$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 369 (pushes old rows 369..452 down to 370..453,
# and grows the used range to A1:R453).
$ws.Rows.Item(369).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A369").Value = 8
$ws.Range("B369").Value = "Terminal La Palmera de La Serena"
$ws.Range("C369").Value = "Coquimbo"
$ws.Range("D369").Value = 44943
$ws.Range("E369").Value = 4
$ws.Range("F369").Value = 100112032
$ws.Range("G369").Value = "Zapallo italiano"
$ws.Range("H369").Value = "Sin especificar"
$ws.Range("I369").Value = "Primera"
$ws.Range("J369").Value = 400
$ws.Range("K369").Value = 7000
$ws.Range("L369").Value = 8000
$ws.Range("M369").Value = 7500
$ws.Range("N369").Value = "$/caja 70 unidades"
$ws.Range("O369").Value = "Provincia del Elquí"
$ws.Range("P369").Value = 107
$ws.Range("Q369").Value = 70
$ws.Range("R369").Value = "Hortaliza"
